$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the K2:L2 values (previously 30000 and 200)
$ws.Range("K2:L2").ClearContents()

# Move the active selection to L2 (was K3)
$ws.Range("L2").Select()
